$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 36 so the existing rows 36-44 shift down to 37-45,
# making room for the new "FY_4.png" record.
$ws.Rows.Item(36).Insert()

# Full target dataset for rows 1-45, columns A-G.
$data = @(
    @("AK_1.png",0.5649999999999999,0.002,0.981,0.429,"Tidak Diketahui","Salah"),
    @("AK_2.png",0.892,0.003,0.993,0.571,"Akhlak Kamiswara","Benar"),
    @("AK_3.png",0.865,0.003,0.985,0.857,"Akhlak Kamiswara","Benar"),
    @("AK_4.png",0.547,0.002,0.986,0.571,"Muhammad Iqbal Baqi","Salah"),
    @("AK_5.png",0.541,0.002,0.992,0.571,"Akhlak Kamiswara","Benar"),
    @("MIB_1.png",1.253,0.004,0.967,0.571,"Muhammad Iqbal Baqi","Benar"),
    @("MIB_2.png",1.114,0.004,0.98,0.429,"Tidak Diketahui","Salah"),
    @("MIB_3.png",1.338,0.004,1,1,"Muhammad Iqbal Baqi","Benar"),
    @("MIB_4.png",1.11,0.004,0.99,0.714,"Muhammad Iqbal Baqi","Benar"),
    @("MIB_5.png",1.246,0.004,0.985,0.571,"Muhammad Iqbal Baqi","Benar"),
    @("AAH_1.png",0.856,0.003,0.978,0.714,"Andrea Ayunove Hutami","Benar"),
    @("AAH_2.png",1.085,0.004,0.994,0.857,"Andrea Ayunove Hutami","Benar"),
    @("AAH_3.png",0.865,0.003,0.978,0.714,"Andrea Ayunove Hutami","Benar"),
    @("TI_1.png",0.95,0.003,0.966,0.429,"Tidak Diketahui","Salah"),
    @("TI_2.png",0.9419999999999999,0.003,0.995,0.429,"Tidak Diketahui","Salah"),
    @("TI_3.png",0.754,0.003,0.994,1,"Toni Ismail","Benar"),
    @("TI_4.png",0.747,0.002,0.99,1,"Toni Ismail","Benar"),
    @("TI_5.png",1.052,0.003,0.994,0.714,"Toni Ismail","Benar"),
    @("RAS_1.png",0.6840000000000001,0.002,0.984,0.429,"Tidak Diketahui","Salah"),
    @("RAS_2.png",1.009,0.003,0.973,1,"Ridha Ayu Salsabila","Benar"),
    @("RAS_3.png",0.572,0.002,0.987,0.571,"Muhammad Iqbal Baqi","Salah"),
    @("RAS_4.png",1.175,0.004,0.973,0.286,"Tidak Diketahui","Salah"),
    @("RAS_5.png",1.104,0.004,0.986,0.857,"Ridha Ayu Salsabila","Benar"),
    @("RR_1.png",1.2,0.004,0.99,0.857,"Rafiqo Rapitasari","Benar"),
    @("RR_2.png",1.173,0.004,0.997,0.286,"Tidak Diketahui","Salah"),
    @("RR_3.png",0.992,0.003,0.906,0.714,"Toni Ismail","Salah"),
    @("RR_4.png",1.199,0.004,0.995,0.571,"Rafiqo Rapitasari","Benar"),
    @("RR_5.png",1.227,0.004,0.993,1,"Rafiqo Rapitasari","Benar"),
    @("AR_1.png",0.884,0.003,0.976,1,"Arizli Romadhon","Benar"),
    @("GA_1.png",1.197,0.004,0.987,1,"Gege Ardiyansyah","Benar"),
    @("GA_2.png",0.678,0.002,0.99,0.571,"Muhammad Iqbal Baqi","Salah"),
    @("GA_3.png",0.767,0.003,0.995,0.429,"Tidak Diketahui","Salah"),
    @("FY_1.png",1.051,0.003,0.98,0.429,"Tidak Diketahui","Salah"),
    @("FY_2.png",1.252,0.004,0.966,0.286,"Tidak Diketahui","Salah"),
    @("FY_3.png",1.278,0.004,0.995,0.429,"Tidak Diketahui","Salah"),
    @("FY_4.png",1.21,0.004,0.985,0.429,"Tidak Diketahui","Salah"),
    @("TO_1.png",0.821,0.003,0.978,0.286,"Tidak Diketahui","Salah"),
    @("TO_2.png",0.9379999999999999,0.003,0.993,0.571,"Tiara Oktavian","Benar"),
    @("TO_3.png",0.876,0.003,0.992,0.714,"Tiara Oktavian","Benar"),
    @("TO_4.png",2.304,0.007,0.853,0.429,"Tidak Diketahui","Salah"),
    @("TO_5.png",2.173,0.007,0.846,0.429,"Tidak Diketahui","Salah"),
    @("TD_1.png",1.545,0.005,0.854,0.286,"Tidak Diketahui","Benar"),
    @("TD_2.png",1.597,0.005,0.772,0.571,"Muhammad Iqbal Baqi","Salah"),
    @("TD_3.png",1.064,0.004,0.981,0.286,"Tidak Diketahui","Benar"),
    @("TD_4.png",0.997,0.003,0.991,0.429,"Tidak Diketahui","Benar")

)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 1
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}
